$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column header "Тип" -> "Молекула" (also updates the Таблица1 tableColumn name)
$ws.Range("E1").Value = "Молекула"

# Replace the text "rad"/"mol" markers in column E with boolean TRUE/FALSE values.
# Radicals (names ending with *) -> FALSE, molecules -> TRUE.
$radicalRows = @(2, 4, 9, 12, 16)
for ($r = 2; $r -le 22; $r++) {
    if ($radicalRows -contains $r) {
        $ws.Range("E$r").Formula = "=FALSE"
    } else {
        $ws.Range("E$r").Formula = "=TRUE"
    }
}
